# "fixed duration contracts implemented"
#
# 1. FD_prices: the fixed-duration contract price column (D2:D73) is moved
#    off its placeholder 0 and populated with the real price (5).
# 2. FD_limits: nothing data-wise changes here, just leave the view settled
#    back at the top of the sheet instead of the scrolled position.
# 3. RawMaterialPrices: the now-unused "forced fd"/"free" scratch
#    calculation (rows 15-16 + the helper formula in row 18) is removed,
#    and the January raw-material price is corrected from 10 to 0.

$wb = $excel.ActiveWorkbook

# --- FD_prices: fill in the fixed-duration contract price for every month ---
$wsFdPrices = $wb.Worksheets.Item("FD_prices")
$wsFdPrices.Activate()
$wsFdPrices.Range("D2:D73").Value = 5
$wsFdPrices.Range("D2:D73").Select()

# --- FD_limits: no data changes, just resettle the view at the top ---
$wsFdLimits = $wb.Worksheets.Item("FD_limits")
$wsFdLimits.Activate()
$wsFdLimits.Range("D19").Select()

# --- RawMaterialPrices: drop the forced-fd/free scratch area, fix C2 ---
$wsRMP = $wb.Worksheets.Item("RawMaterialPrices")
$wsRMP.Activate()
$wsRMP.Range("C2").Value = 0
$wsRMP.Range("E15:F16").ClearContents()
$wsRMP.Range("E18").ClearContents()
$wsRMP.Range("F17").Select()
